$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

$ws.Range("E2").Value = "2026-02-12 17:57:16"
$ws.Range("E3").Value = "2026-02-12 17:57:19"
$ws.Range("E4").Value = "2026-02-12 17:57:21"
$ws.Range("J4").Value = "998.3 hPa"
$ws.Range("O4").Value = "16.4 °C"
$ws.Range("E5").Value = "2026-02-12 17:57:24"
$ws.Range("E6").Value = "2026-02-12 17:57:26"
$ws.Range("J6").Value = "998.1 hPa"
$ws.Range("E7").Value = "2026-02-12 17:57:28"
$ws.Range("J7").Value = "1001.1 hPa"
$ws.Range("K7").Value = "14.1 MJ/m2"
$ws.Range("E8").Value = "2026-02-12 17:57:31"
$ws.Range("J8").Value = "1000.3 hPa"
$ws.Range("E9").Value = "2026-02-12 17:57:33"
$ws.Range("H9").Value = "'70%"
$ws.Range("E10").Value = "2026-02-12 17:57:36"
$ws.Range("E11").Value = "2026-02-12 17:57:38"
$ws.Range("H11").Value = "'49%"
$ws.Range("O11").Value = "9.6 °C"
$ws.Range("E12").Value = "2026-02-12 17:57:41"
$ws.Range("H12").Value = "'76%"
$ws.Range("O12").Value = "12.6 °C"
$ws.Range("E13").Value = "2026-02-12 17:57:43"
$ws.Range("J13").Value = "1000.9 hPa"
$ws.Range("E14").Value = "2026-02-12 17:57:45"
$ws.Range("E15").Value = "2026-02-12 17:57:48"
$ws.Range("H15").Value = "'57%"
$ws.Range("O15").Value = "13.9 °C"
$ws.Range("E16").Value = "2026-02-12 17:57:50"
$ws.Range("O16").Value = "-4.6 °C"
$ws.Range("E17").Value = "2026-02-12 17:57:53"
$ws.Range("E18").Value = "2026-02-12 17:57:55"
$ws.Range("J18").Value = "998.6 hPa"
$ws.Range("K18").Value = "13.9 MJ/m2"
$ws.Range("E19").Value = "2026-02-12 17:57:58"
$ws.Range("H19").Value = "'63%"
$ws.Range("E20").Value = "2026-02-12 17:58:00"
$ws.Range("E21").Value = "2026-02-12 17:58:02"
$ws.Range("J21").Value = "1001.4 hPa"
$ws.Range("E22").Value = "2026-02-12 17:58:05"
$ws.Range("E23").Value = "2026-02-12 17:58:07"
$ws.Range("K23").Value = "10.9 MJ/m2"
$ws.Range("E24").Value = "2026-02-12 17:58:10"
$ws.Range("J24").Value = "1006.1 hPa"
$ws.Range("O24").Value = "11.7 °C"
$ws.Range("E25").Value = "2026-02-12 17:58:12"
$ws.Range("E26").Value = "2026-02-12 17:58:15"
$ws.Range("J26").Value = "997.5 hPa"
$ws.Range("O26").Value = "6.2 °C"
$ws.Range("E27").Value = "2026-02-12 17:58:18"
$ws.Range("E28").Value = "2026-02-12 17:58:20"
$ws.Range("J28").Value = "997.8 hPa"
$ws.Range("E29").Value = "2026-02-12 17:58:23"
$ws.Range("E30").Value = "2026-02-12 17:58:25"
$ws.Range("J30").Value = "998.3 hPa"
$ws.Range("E31").Value = "2026-02-12 17:58:28"
$ws.Range("J31").Value = "997.8 hPa"
$ws.Range("E32").Value = "2026-02-12 17:58:30"
$ws.Range("K32").Value = "14.0 MJ/m2"
$ws.Range("E33").Value = "2026-02-12 17:58:33"
$ws.Range("H33").Value = "'53%"
$ws.Range("J33").Value = "1000.6 hPa"
$ws.Range("O33").Value = "6.8 °C"
$ws.Range("E34").Value = "2026-02-12 17:58:35"
$ws.Range("E35").Value = "2026-02-12 17:58:38"
$ws.Range("J35").Value = "1007.5 hPa"
$ws.Range("E36").Value = "2026-02-12 17:58:40"
$ws.Range("H36").Value = "'64%"
$ws.Range("J36").Value = "998.7 hPa"
$ws.Range("O36").Value = "14.2 °C"
$ws.Range("E37").Value = "2026-02-12 17:58:43"
$ws.Range("J37").Value = "999.0 hPa"
$ws.Range("O37").Value = "10.5 °C"
$ws.Range("E38").Value = "2026-02-12 17:58:45"
$ws.Range("E39").Value = "2026-02-12 17:58:48"
$ws.Range("E40").Value = "2026-02-12 17:58:50"
$ws.Range("J40").Value = "1002.1 hPa"
$ws.Range("O40").Value = "10.1 °C"
$ws.Range("E41").Value = "2026-02-12 17:58:53"
$ws.Range("J41").Value = "1005.1 hPa"
$ws.Range("E42").Value = "2026-02-12 17:58:55"
$ws.Range("E43").Value = "2026-02-12 17:58:58"
$ws.Range("E44").Value = "2026-02-12 17:59:00"
$ws.Range("H44").Value = "'68%"
$ws.Range("I44").Value = "0.4 mm"
$ws.Range("O44").Value = "-3.0 °C"
$ws.Range("E45").Value = "2026-02-12 17:59:03"
$ws.Range("J45").Value = "1004.2 hPa"
$ws.Range("E46").Value = "2026-02-12 17:59:05"
$ws.Range("J46").Value = "1006.9 hPa"
